# Apply the "28-5-2024 create category list done" edit to myWorkRecord.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in row 21 (2024-05-28 entry) which previously only had sr no / date.
# This must happen BEFORE the row-20 typo fix so the shared-string table
# ends up with the same ordering as the authored workbook.
$ws.Range("C21").Value = "make category component and make routing and make table with search bar and add button"
$ws.Range("D21").Value = "no"
$ws.Range("E21").Value = "1 day"

# Fix the typo in row 20 (2024-05-27 entry): "changes butto text color" -> "changes button text color"
$ws.Range("C20").Value = "changes button text color"

# Add a new row 22 (2024-05-29 entry) with sr no and date only
$ws.Range("A22").Value = 21
$ws.Range("B21").Copy()
$ws.Range("B22").PasteSpecial(-4122)  # xlPasteFormats - reuse the existing date style (s=1)
$ws.Range("B22").Value = 45441
$excel.CutCopyMode = $false

# Update the view so the newly active cell / scroll position matches
$ws.Range("C22").Select()
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 3
